# Add bottom border to the 'it will create the user table...' paragraph
$d = $word.ActiveDocument
$pTable = $d.Paragraphs(2)
$pTable.Borders(-3).LineStyle = 1
$pTable.Borders(-3).LineWidth = 3
$pTable.Borders(-3).ColorIndex = 0
$pTable.Borders.DistanceFromBottom = 1

# Make sure we have enough empty paragraphs after it to host the new content.
# There are currently 5 empty paragraphs (3..7); we need 5 content paragraphs
# followed by 3 empty ones, i.e. 3 more paragraphs must be added.
$anchor = $d.Paragraphs(4)
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()

$xml0 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>JWT tokens and SMTP :</w:t>
      </w:r>
    </w:p>
'@
$d.Paragraphs(3).Range.InsertXML($xml0)

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Json web tokens  : user id as a data and the SH256 algorithm and  </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>32 bit</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> hexadecimal secret key to encode the data using JWT . using this we create the token and send the user registered email using </w:t>
      </w:r>
    </w:p>
'@
$d.Paragraphs(4).Range.InsertXML($xml1)

$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> simple mail transfer protocol we send the mail with verification  link is link the (http:127.0.0.1:8000//user.verify/?token=”asdfafa65a46dfa5s6df4a5fd6a5dfa6fd5a4asd”)  </w:t>
      </w:r>
    </w:p>
'@
$d.Paragraphs(5).Range.InsertXML($xml2)

$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">for secret key : </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t xml:space="preserve">import secrets </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t xml:space="preserve">    &gt; </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>secrets.token_hex</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>(32)</w:t>
      </w:r>
    </w:p>
'@
$d.Paragraphs(6).Range.InsertXML($xml3)

$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pBdr>
          <w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/>
        </w:pBdr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">after sending the mail to the user . then user can hit that link then the verify </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>api</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> is automatically called . in the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">verify </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>api</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> we make the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>is</w:t>
      </w:r>
      <w:r>
        <w:t>_</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>verified</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">filled </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>In</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> the database as the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>Trur</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
'@
$d.Paragraphs(7).Range.InsertXML($xml4)

